$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "quantity" values for the changed rows (row -> new value)
$updates = @{
    2  = 205591
    4  = 91837
    5  = 81729
    8  = 64768
    9  = 54652
    10 = 48417
    12 = 40107
    16 = 27827
    18 = 24527
    19 = 23156
    20 = 19630
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# Recompute the "periodicity" column (C) as quantity / sum(quantity) for all data rows (2..35)
$total = 0
for ($r = 2; $r -le 35; $r++) {
    $total += $ws.Cells.Item($r, 2).Value()
}

for ($r = 2; $r -le 35; $r++) {
    $q = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 3).Value = $q / $total
}
